# Refresh the cryptos price/volume snapshot (Price = col D, Volume(1h) = col E).
# A leading apostrophe forces Excel to keep numeric-looking Price strings
# (e.g. "1.00", "0.508") as literal text instead of coercing them to numbers,
# matching how the source feed writes them (fixed-format strings, not numerics).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.120.32'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '2.432.54'
$ws.Range("E3").Value = '  +7.27%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''296.14'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = '''95.80'
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("E7").Value = '  +1.50%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.508'
$ws.Range("D10").Value = '''35.35'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("D12").Value = '''7.13'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").Value = '2.803.67'
$ws.Range("E14").Value = '  +7.30%  '
$ws.Range("D15").Value = '2.435.11'
$ws.Range("E15").Value = '  +7.40%  '
$ws.Range("E16").Value = '  +6.42%  '
$ws.Range("D17").Value = '''14.17'
$ws.Range("E17").Value = '  +3.92%  '
$ws.Range("D18").Value = '46.005.50'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '''12.53'
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("D20").Value = '0.0₃0947'
$ws.Range("E20").Value = '  -1.76%  '
$ws.Range("D21").Value = '''6.21'
$ws.Range("E21").Value = '  +6.94%  '
$ws.Range("D22").Value = '''67.42'
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").Value = '''244.67'
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("D24").Value = '''2.80'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  +4.94%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '''39.60'
$ws.Range("E27").Value = '  -5.00%  '
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("D29").Value = '''9.76'
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").Value = '''3.88'
$ws.Range("E30").Value = '  +17.36%  '
$ws.Range("D31").Value = '''21.30'
$ws.Range("E31").Value = '  +5.23%  '
$ws.Range("D32").Value = '''2.78'
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = '''148.33'
$ws.Range("E33").Value = '  +1.33%  '
$ws.Range("E34").Value = '  +2.71%  '
$ws.Range("D35").Value = '''0.0771'
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("D36").Value = '''2.02'
$ws.Range("E36").Value = '  +19.48%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").Value = '''14.83'
$ws.Range("E39").Value = '  -5.88%  '
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = '''0.0301'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("D42").Value = '''3.26'
$ws.Range("E42").Value = '  +4.56%  '
$ws.Range("D43").Value = '1.981.82'
$ws.Range("E43").Value = '  +11.24%  '
$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''90.03'
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").Value = '''1.81'
$ws.Range("E46").Value = '  -4.77%  '
$ws.Range("D47").Value = '''16.45'
$ws.Range("E47").Value = '  +29.09%  '
$ws.Range("D48").Value = '''8.63'
$ws.Range("E48").Value = '  +10.23%  '
$ws.Range("D49").Value = '''101.23'
$ws.Range("E49").Value = '  +7.45%  '
$ws.Range("D50").Value = '2.672.85'
$ws.Range("E50").Value = '  +7.28%  '
$ws.Range("E51").Value = '  +0.07%  '
